$wb = $excel.ActiveWorkbook

# --- Populate the "Quiz 1" sheet with grade data ---
$quiz = $wb.Worksheets.Item("Quiz 1")

$quiz.Range("A1").Value = "John Doe"
$quiz.Range("B1").Value = "'10/10"
$quiz.Range("B1").NumberFormat = "d-mmm"

$quiz.Range("A2").Value = "Molly Doe"
$quiz.Range("B2").Value = 98.5

$quiz.Range("A3").Value = "Stephen Jane"
$quiz.Range("B3").Value = 45

# --- Size columns to fit the new content (best-fit widths) ---
$quiz.Columns.Item(1).ColumnWidth = 12
$quiz.Columns.Item(2).ColumnWidth = 5

# --- Update selections on the other sheets ---
$roster = $wb.Worksheets.Item("Roster")
$roster.Range("B4").Select() | Out-Null

# --- Make "Quiz 1" the active sheet/tab and set its selection ---
$quiz.Activate() | Out-Null
$quiz.Range("A4").Select() | Out-Null
